$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hárok1")

# Fill in new row 10 data (test case "4 to 1")
$ws.Range("A10").Value = "4 to 1"
$ws.Range("C10").Value = "ERR"
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = "Na dvojek zastavilo s hlaskou ERR_1SPOM_C"

# Move the active selection to E11 (next empty description cell)
$ws.Range("E11").Select()
